# Automatic update of files.
# Row 2 and Row 3 in the "Artfynd" sheet effectively swap their species
# records (columns A, D, E, F, G, H, I, K, Q, R), while the
# "Taxonsorteringsordning" column (B) is updated to new, independent
# values (not simply swapped).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: becomes the "Grönpyrola" / Pyrola chlorantha record ---
$ws.Range("A2").Value = 110915106
$ws.Range("B2").Value = 103288
$ws.Range("D2").Value = "LC"
$ws.Range("E2").Value = 221144
$ws.Range("F2").Value = "Grönpyrola"
$ws.Range("G2").Value = "Pyrola chlorantha"
$ws.Range("H2").Value = "Sw."

$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "50"
$ws.Range("I2").ClearFormats()

$ws.Range("K2").Value = "blomning"
$ws.Range("Q2").Value = 585461.8925125685
$ws.Range("R2").Value = 6315287.846391106

# --- Row 3: becomes the "Knärot" / Goodyera repens record ---
$ws.Range("A3").Value = 110915063
$ws.Range("B3").Value = 96348
$ws.Range("D3").Value = "VU"
$ws.Range("E3").Value = 220787
$ws.Range("F3").Value = "Knärot"
$ws.Range("G3").Value = "Goodyera repens"
$ws.Range("H3").Value = "(L.) R. Br."

$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = "8"
$ws.Range("I3").ClearFormats()

$ws.Range("K3").Value = "fullt utvecklade blad"
$ws.Range("Q3").Value = 585441.8752236688
$ws.Range("R3").Value = 6315225.284495098
